# Generate Report for Handoff
# Updates the localization-status workbook to reflect that the
# 646ebfcf-678c-4f34-baf2-3860a32c3582.md and 81d86a05-ef97-45a1-ae7e-6b1389d4d1ee.md
# files are "Ready for handoff", with refreshed handoff timestamps and an
# error detail noting the handback file version is stale.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# ---------------------------------------------------------------------
# Overview sheet: rows 4 (646ebfcf...) and 5 (81d86a05...)
#   E = zh-cn column, F = de-de column, G = Latest HO Xliff Generate Date
# ---------------------------------------------------------------------
$overview.Range("E4").Value = "Ready for handoff"
$overview.Range("F4").Value = "Ready for handoff"
$overview.Range("G4").Value = "2016-10-27 10:24:39"

$overview.Range("E5").Value = "Ready for handoff"
$overview.Range("F5").Value = "Ready for handoff"
$overview.Range("G5").Value = "2016-10-27 10:24:39"

# ---------------------------------------------------------------------
# zh-cn sheet: rows 4 (646ebfcf...) and 5 (81d86a05...)
#   C = Status, H = Latest Handoff Datetime, P = Error Detail
# ---------------------------------------------------------------------
$zhcn.Range("C4").Value = "Ready for handoff"
$zhcn.Range("H4").Value = "2016-10-27 10:24:26"
$zhcn.Range("P4").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f4a4405ebf8498665d144fb33aa1061f72dec67e/e2e/646ebfcf-678c-4f34-baf2-3860a32c3582.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6fb3941edffc46f2380f017ff80a12b20f40a626/e2e/646ebfcf-678c-4f34-baf2-3860a32c3582.md."

$zhcn.Range("C5").Value = "Ready for handoff"
$zhcn.Range("H5").Value = "2016-10-27 10:24:26"
$zhcn.Range("P5").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f4a4405ebf8498665d144fb33aa1061f72dec67e/e2e/81d86a05-ef97-45a1-ae7e-6b1389d4d1ee.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6fb3941edffc46f2380f017ff80a12b20f40a626/e2e/81d86a05-ef97-45a1-ae7e-6b1389d4d1ee.md."

# Widen the Error Detail column so the new message is readable.
$zhcn.Columns.Item(16).ColumnWidth = 39.17

# ---------------------------------------------------------------------
# de-de sheet: rows 4 (646ebfcf...) and 5 (81d86a05...)
#   C = Status, H = Latest Handoff Datetime, P = Error Detail
# ---------------------------------------------------------------------
$dede.Range("C4").Value = "Ready for handoff"
$dede.Range("H4").Value = "2016-10-27 10:24:39"
$dede.Range("P4").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f4a4405ebf8498665d144fb33aa1061f72dec67e/e2e/646ebfcf-678c-4f34-baf2-3860a32c3582.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6fb3941edffc46f2380f017ff80a12b20f40a626/e2e/646ebfcf-678c-4f34-baf2-3860a32c3582.md."

$dede.Range("C5").Value = "Ready for handoff"
$dede.Range("H5").Value = "2016-10-27 10:24:39"
$dede.Range("P5").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f4a4405ebf8498665d144fb33aa1061f72dec67e/e2e/81d86a05-ef97-45a1-ae7e-6b1389d4d1ee.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6fb3941edffc46f2380f017ff80a12b20f40a626/e2e/81d86a05-ef97-45a1-ae7e-6b1389d4d1ee.md."

# Widen the Error Detail column so the new message is readable.
$dede.Columns.Item(16).ColumnWidth = 39.17
